$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row and fix "de/del/el/la/los/las/y" capitalization in municipality names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B5").Value = "Pabellón De Arteaga"
$ws.Range("B6").Value = "Rincón De Romos"
$ws.Range("B21").Value = "Amatenango De La Frontera"
$ws.Range("B30").Value = "Comitán De Domínguez"
$ws.Range("B44").Value = "Ocozocoautla De Espinosa"
$ws.Range("B76").Value = "Hidalgo Del Parral"
$ws.Range("B79").Value = "San Francisco De Borja"
$ws.Range("B93").Value = "San Juan De Sabinas"
$ws.Range("A103").Value = "Ciudad De México"
$ws.Range("B107").Value = "Cuajimalpa De Morelos"
$ws.Range("A128").Value = "Estado De México"
$ws.Range("B128").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B129").Value = "Almoloya De Alquisiras"
$ws.Range("B130").Value = "Atizapán De Zaragoza"
$ws.Range("B135").Value = "Ecatepec De Morelos"
$ws.Range("B137").Value = "Naucalpan De Juárez"
$ws.Range("B149").Value = "Tlalnepantla De Baz"
$ws.Range("B154").Value = "Valle De Bravo"
$ws.Range("B160").Value = "San Miguel De Allende"
$ws.Range("B161").Value = "Apaseo El Alto"
$ws.Range("B162").Value = "Apaseo El Grande"
$ws.Range("B168").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B180").Value = "San Diego De La Unión"
$ws.Range("B182").Value = "San Francisco Del Rincón"
$ws.Range("B184").Value = "San Luis De La Paz"
$ws.Range("B185").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B189").Value = "Valle De Santiago"
$ws.Range("B193").Value = "Acapulco De Juárez"
$ws.Range("B194").Value = "Ajuchitlán Del Progreso"
$ws.Range("B197").Value = "Atoyac De Álvarez"
$ws.Range("B198").Value = "Ayutla De Los Libres"
$ws.Range("B199").Value = "Chilapa De Álvarez"
$ws.Range("B200").Value = "Chilpancingo De Los Bravo"
$ws.Range("B202").Value = "Coyuca De Catalán"
$ws.Range("B205").Value = "Huitzuco De Los Figueroa"
$ws.Range("B206").Value = "Iguala De La Independencia"
$ws.Range("B215").Value = "Taxco De Alarcón"
$ws.Range("B217").Value = "Técpan De Galeana"
$ws.Range("B219").Value = "Tepecoacuilco De Trujano"
$ws.Range("B223").Value = "Tlapa De Comonfort"
$ws.Range("B232").Value = "Atotonilco El Grande"
$ws.Range("B234").Value = "Cuautepec De Hinojosa"
$ws.Range("B236").Value = "Huejutla De Reyes"
$ws.Range("B241").Value = "Mineral Del Chico"
$ws.Range("B242").Value = "Molango De Escamilla"
$ws.Range("B244").Value = "Pachuca De Soto"
$ws.Range("B249").Value = "Tenango De Doria"
$ws.Range("B250").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B252").Value = "Tezontepec De Aldama"
$ws.Range("B256").Value = "Zacualtipán De Ángeles"
$ws.Range("B267").Value = "Jilotlán De Los Dolores"
$ws.Range("B274").Value = "San Juanito De Escobedo"
$ws.Range("B275").Value = "Tamazula De Gordiano"
$ws.Range("B277").Value = "Tepatitlán De Morelos"
$ws.Range("B279").Value = "Tizapán El Alto"
$ws.Range("B334").Value = "Amatlán De Cañas"
$ws.Range("B350").Value = "San Nicolás De Los Garza"
$ws.Range("B353").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B354").Value = "Coicoyán De Las Flores"
$ws.Range("B357").Value = "Eloxochitlán De Flores Magón"
$ws.Range("B358").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B359").Value = "Ixtlán De Juárez"
$ws.Range("B360").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B362").Value = "Mártires De Tacubaya"
$ws.Range("B363").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B364").Value = "Oaxaca De Juárez"
$ws.Range("B366").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B379").Value = "San Mateo Del Mar"
$ws.Range("B380").Value = "San Miguel Del Puerto"
$ws.Range("B402").Value = "Santo Domingo De Morelos"
$ws.Range("B408").Value = "Tanetze De Zaragoza"
$ws.Range("B409").Value = "Tataltepec De Valdés"
$ws.Range("B410").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B411").Value = "Tlacolula De Matamoros"
$ws.Range("B412").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B414").Value = "Zimatlán De Álvarez"
$ws.Range("B423").Value = "Huehuetlán El Grande"
$ws.Range("B425").Value = "Izúcar De Matamoros"
$ws.Range("B426").Value = "Los Reyes De Juárez"
$ws.Range("B429").Value = "Palmar De Bravo"
$ws.Range("B434").Value = "San Salvador El Verde"
$ws.Range("B435").Value = "Tecali De Herrera"
$ws.Range("B448").Value = "Amealco De Bonfil"
$ws.Range("B450").Value = "Cadereyta De Montes"
$ws.Range("B452").Value = "Jalpan De Serra"
$ws.Range("B453").Value = "Pinal De Amoles"
$ws.Range("B456").Value = "San Juan Del Río"
$ws.Range("B463").Value = "Axtla De Terrazas"
$ws.Range("B465").Value = "Cerro De San Pedro"
$ws.Range("B466").Value = "Ciudad Del Maíz"
$ws.Range("B474").Value = "Mexquitic De Carmona"
$ws.Range("B478").Value = "San Ciro De Acosta"
$ws.Range("B483").Value = "Santa María Del Río"
$ws.Range("B484").Value = "Soledad De Graciano Sánchez"
$ws.Range("B491").Value = "Villa De Arriaga"
$ws.Range("B492").Value = "Villa De Ramos"
$ws.Range("B508").Value = "Nacozari De García"
$ws.Range("B540").Value = "Soto La Marina"
$ws.Range("B549").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B558").Value = "Amatlán De Los Reyes"
$ws.Range("B565").Value = "Castillo De Teayo"
$ws.Range("B572").Value = "Cosamaloapan De Carpio"
$ws.Range("B581").Value = "Hueyapan De Ocampo"
$ws.Range("B582").Value = "Ignacio De La Llave"
$ws.Range("B585").Value = "Ixhuatlán Del Café"
$ws.Range("B593").Value = "Martínez De La Torre"
$ws.Range("B602").Value = "Paso Del Macho"
$ws.Range("B604").Value = "Poza Rica De Hidalgo"
$ws.Range("B609").Value = "Sayula De Alemán"
$ws.Range("B611").Value = "Soledad De Doblado"
$ws.Range("B614").Value = "Tatahuicapan De Juárez"
$ws.Range("B631").Value = "Vega De Alatorre"
$ws.Range("B650").Value = "Villa De Cos"

# Remove trailing metadata/footer rows (655:659) that are no longer part of the data table
$ws.Rows("655:659").Delete()

